$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 'Test zur Rhythmuswahrnehmung in der Sprache'

$ws.Cells.Item(3, 2).Value = 'Wir werden Ihnen nun einige Reihen von jeweils vier einzelnen Wörtern oder kurzen Ausdrücken präsentieren, z.B. Teller – Schachtel – Ratte – Parkett. Bitte lesen Sie diese Reihen aufmerksam und achten Sie dabei besonders auf den Sprachrhythmus, d.h. auf die Betonungsmuster der Wörter. <br/> In manchen dieser Reihen weisen alle Elemente, d.h. alle enthaltenen Einzelwörter oder Mehrwort-Ausdrücke, dasselbe Betonungsmuster auf.  In einigen dieser Reihen ist jedoch genau ein Element enthalten, das ein anderes Betonungsmuster aufweist als die übrigen – im Beispiel oben ist es das  Wort „Parkett“. Ihre Aufgabe besteht darin, für jede Reihe zu beurteilen, ob ein Element mit abweichendem Betonungsmuster enthalten ist oder nicht. Falls alle Elemente der Reihe das gleiche Betonungsmuster aufweisen, klicken Sie bitte auf das Feld „Alle gleich“. Falls ein Element ein abweichendes Betonungsmuster aufweist, klicken Sie bitte auf das Feld dieses abweichenden Elements.<br/> Bitte beachten Sie, dass Sie für jede Reihe **10 Sekunden Zeit** haben, um Ihr Urteil abzugeben.'
$ws.Cells.Item(3, 3).Value = 'We will now present you with some series of four individual words or short expressions each, e.g. plate - box - rat - parquet.  Please read these series carefully, paying particular  attention to the rhythm of speech, i.e. the stress patterns of the words. <br/> In some of these series, however, there is exactly one element that has a different stress pattern than the others - in the example above it is the word "parquet".  Your task is to judge for each row whether or not there is an element with a different stress pattern.  If all the elements in the row have the same stress pattern, please click on the field "All equal".  If an element has a different stress pattern, please click on the field of this different element. <br/>Please note that you have **10 seconds** for each row to make your judgement.'

$ws.Cells.Item(4, 2).Value = 'Ist der Rhythmus aller Wörter gleich oder hat eines ein abweichendes Betonungsmuster?<br/>  Wählen Sie “Alle gleich” oder klicken Sie auf das Wort, das abweicht. <br/> Sie haben {{time_out}} Sekunden Zeit zu antworten.' + "`n"

$ws.Cells.Item(9, 2).Value = 'Test: Rhythmuswahrnehmung in der Sprache!'

$ws.Cells.Item(14, 2).Value = 'Korrekt, ReVIER hat ein abweichendes Betonungsmuster.'

# Row 4 formatting: taller row, wrap text + top vertical alignment for B4:D4
$ws.Rows.Item(4).RowHeight = 96
$ws.Range("B4:D4").WrapText = $true
$ws.Range("B4:D4").VerticalAlignment = -4160

# New column C width
$ws.Columns.Item(3).ColumnWidth = 32.25

# Update selected cell
$ws.Range("B9").Select()
